$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Equity Typel" column (column B) entirely, shifting remaining
# columns (Endorsement, Endorsement Folio, ... Company Master Id*) left.
$ws.Range("B1").EntireColumn.Delete()
